# CIV-6625 Update GA order template
# Remove the "Classification: Controlled" text-box shape from the
# primary (default) footer.

$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1 -> this is the "default" footer (footer2.xml),
# the one that originally carried the classification text box.
$ftr = $d.Sections(1).Footers(1)

for ($i = $ftr.Shapes.Count; $i -ge 1; $i--) {
    $shp = $ftr.Shapes($i)
    if ($shp.Name -eq "Text Box 4" -or $shp.AlternativeText -eq "Classification: Controlled") {
        $shp.Delete()
    }
}
